$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Control 0)
$ws.Range("D2").Value = 7.654292115652291 / [Math]::Pow(10, 9)
$ws.Range("E2").Value = 7.654292115652291 / [Math]::Pow(10, 9)

# Row 3 (Control 6)
$ws.Range("C3").Value = $false
$ws.Range("D3").Value = 0.9981312601615195
$ws.Range("E3").Value = 0.9981312601615195

# Row 4 (Control 9)
$ws.Range("D4").Value = 4.621211096430835 / [Math]::Pow(10, 9)
$ws.Range("E4").Value = 4.621211096430835 / [Math]::Pow(10, 9)

# Row 5 (Control 24)
$ws.Range("D5").Value = 0.999999998714433
$ws.Range("E5").Value = 0.999999998714433

# Row 6 (Control 32)
$ws.Range("D6").Value = 0.9999999999999925
$ws.Range("E6").Value = 0.9999999999999925

# Row 7 (MDD 36)
$ws.Range("D7").Value = 0.9999999999999443
$ws.Range("E7").Value = 5.573319583618286 / [Math]::Pow(10, 14)

# Row 8 (MDD 10)
$ws.Range("D8").Value = 0.9999999991977113
$ws.Range("E8").Value = 8.022886799352591 / [Math]::Pow(10, 10)

# Row 9 (MDD 39)
$ws.Range("D9").Value = 0.9999999999930884
$ws.Range("E9").Value = 6.91158241750145 / [Math]::Pow(10, 12)

# Row 10 (MDD 14)
$ws.Range("D10").Value = 0.9999999998388487
$ws.Range("E10").Value = 1.611513145149956 / [Math]::Pow(10, 10)

# Row 11 (MDD 18)
$ws.Range("F11").Value = 5.927690029144287
$ws.Range("G11").Value = 0.7
